$wb = $excel.ActiveWorkbook

# ---- Schedule sheet updates (rows 2-4 reshuffled after re-optimisation run) ----
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("A2").Value2 = 46067.04166666666
$wsSchedule.Range("B2").Value2 = 46067.22916666666
$wsSchedule.Range("C2").Value2 = 4.5
$wsSchedule.Range("D2").Value2 = 17.01
$wsSchedule.Range("E2").Value2 = 499.4363595
$wsSchedule.Range("F2").Value2 = 29.36133800705468
$wsSchedule.Range("A3").Value2 = 46067.3125
$wsSchedule.Range("B3").Value2 = 46067.79166666666
$wsSchedule.Range("C3").Value2 = 11.5
$wsSchedule.Range("D3").Value2 = 43.47
$wsSchedule.Range("E3").Value2 = 681.1232219999997
$wsSchedule.Range("F3").Value2 = 15.66881118012422
$wsSchedule.Range("A4").Value2 = 46068.27083333334
$wsSchedule.Range("C4").Value2 = 12
$wsSchedule.Range("D4").Value2 = 45.36
$wsSchedule.Range("E4").Value2 = 462.74472075
$wsSchedule.Range("F4").Value2 = 10.20160319113757

# ---- Detailed sheet updates (re-run 176 values) ----
$wsDetailed = $wb.Worksheets.Item("Detailed")
$wsDetailed.Range("E4").Value = 'ON'
$wsDetailed.Range("E5").Value = 'ON'
$wsDetailed.Range("E6").Value = 'ON'
$wsDetailed.Range("E7").Value = 'ON'
$wsDetailed.Range("E8").Value = 'ON'
$wsDetailed.Range("E9").Value = 'ON'
$wsDetailed.Range("E10").Value = 'ON'
$wsDetailed.Range("E11").Value = 'ON'
$wsDetailed.Range("E12").Value = 'ON'
$wsDetailed.Range("B37").Value2 = 57.09
$wsDetailed.Range("B38").Value2 = 57.09
$wsDetailed.Range("B39").Value2 = 56.98
$wsDetailed.Range("C39").Value = 'historical'
$wsDetailed.Range("B40").Value2 = 63.60109
$wsDetailed.Range("C40").Value = 'historical'
$wsDetailed.Range("E40").Value = 'OFF'
$wsDetailed.Range("B41").Value2 = 64.89
$wsDetailed.Range("C41").Value = 'historical'
$wsDetailed.Range("B42").Value2 = 64.89
$wsDetailed.Range("C42").Value = 'historical'
$wsDetailed.Range("E42").Value = 'OFF'
$wsDetailed.Range("B43").Value2 = 65.00005
$wsDetailed.Range("C43").Value = 'historical'
$wsDetailed.Range("E43").Value = 'OFF'
$wsDetailed.Range("B44").Value2 = 64.89
$wsDetailed.Range("C44").Value = 'historical'
$wsDetailed.Range("E44").Value = 'OFF'
$wsDetailed.Range("B45").Value2 = 57.60478
$wsDetailed.Range("C45").Value = 'historical'
$wsDetailed.Range("E45").Value = 'OFF'
$wsDetailed.Range("C46").Value = 'historical'
$wsDetailed.Range("E46").Value = 'OFF'
$wsDetailed.Range("B47").Value2 = 57.09
$wsDetailed.Range("C47").Value = 'historical'
$wsDetailed.Range("E47").Value = 'OFF'
$wsDetailed.Range("B48").Value2 = 56.98
$wsDetailed.Range("C48").Value = 'historical'
$wsDetailed.Range("E48").Value = 'OFF'
$wsDetailed.Range("B49").Value2 = 56.98
$wsDetailed.Range("E49").Value = 'OFF'
$wsDetailed.Range("B50").Value2 = 53.12933
$wsDetailed.Range("E50").Value = 'OFF'
$wsDetailed.Range("B52").Value2 = 53.53295
$wsDetailed.Range("B53").Value2 = 53.94375
$wsDetailed.Range("B54").Value2 = 54.01585
$wsDetailed.Range("B55").Value2 = 52.6617
$wsDetailed.Range("B56").Value2 = 52.37068
$wsDetailed.Range("B57").Value2 = 52.37048
$wsDetailed.Range("B58").Value2 = 54.43034
$wsDetailed.Range("B59").Value2 = 52.0279
$wsDetailed.Range("B60").Value2 = 53.70872
$wsDetailed.Range("B61").Value2 = 56.98
$wsDetailed.Range("B62").Value2 = 56.03
$wsDetailed.Range("B63").Value2 = 53.68323
$wsDetailed.Range("E63").Value = 'ON'
$wsDetailed.Range("B65").Value2 = 36.0601
$wsDetailed.Range("B72").Value2 = 22.07
$wsDetailed.Range("B73").Value2 = 22.07
$wsDetailed.Range("B74").Value2 = 1.23444
$wsDetailed.Range("B75").Value2 = 0.03385
$wsDetailed.Range("B76").Value2 = 0.51
$wsDetailed.Range("B77").Value2 = 0.36345
$wsDetailed.Range("B78").Value2 = 0.50984
$wsDetailed.Range("B79").Value2 = -6.99601
$wsDetailed.Range("B80").Value2 = -1.80118
$wsDetailed.Range("B81").Value2 = 0.0103
$wsDetailed.Range("B82").Value2 = 9.023770000000001
$wsDetailed.Range("B83").Value2 = 8.768990000000001
$wsDetailed.Range("B84").Value2 = 8.92048
$wsDetailed.Range("B85").Value2 = 31.10246
$wsDetailed.Range("B86").Value2 = 36.62555
$wsDetailed.Range("B87").Value2 = 57.6508
$wsDetailed.Range("B89").Value2 = 68.62712999999999
$wsDetailed.Range("B90").Value2 = 68.85534
$wsDetailed.Range("B92").Value2 = 64.89
$wsDetailed.Range("B93").Value2 = 60.64516
$wsDetailed.Range("B94").Value2 = 56.98
$wsDetailed.Range("B95").Value2 = 47.60427
$wsDetailed.Range("B96").Value2 = 42.12537
